$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = "heba@med.asu.edu.eg, servinaz@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm, Amira.Sobhy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg"
$ws.Range("G3").Value = "backup@backdoor.com, alshimaa.atef@med.asu.edu.egm, manar.montaser@med.asu.edu.eg, gehanadel@med.asu.edu.eg"
$ws.Range("G4").Value = "heba@med.asu.edu.eg, hananragab@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg"
$ws.Range("G5").Value = "AbeerRagheb@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg"
$ws.Range("G9").Value = "Marina.youhana@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg"
$ws.Range("G10").Value = "amira.m.ibrahim@med.asu.edu.eg, basma.hamed@med.asu.edu.eg"
$ws.Range("G12").Value = "dalia.t.abdelmegid@med.asu.edu.eg, heba.e.mohamed@med.asu.edu.eg, mona.I.hussein@med.asu.edu.eg"
$ws.Range("G13").Value = "youstina.gamil@med.asu.edu.eg, Sarah.Mahdy@med.asu.edu.eg"
$ws.Range("G17").Value = "enas.omran@med.asu.edu.eg, user@user.com, Walaa.h.ghanima@med.asu.edu.eg"
$ws.Range("G18").Value = "abdullah.elagrody@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, eman.samir@med.asu.edu.eg"
$ws.Range("G19").Value = "marinasorial@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, eman.samir@med.asu.edu.eg"
$ws.Range("G20").Value = "marinasorial@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, marina_atef@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg"
$ws.Range("G21").Value = "heba@med.asu.edu.eg, servinaz@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm, Amira.Sobhy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg"
$ws.Range("G22").Value = "backup@backdoor.com, alshimaa.atef@med.asu.edu.egm, manar.montaser@med.asu.edu.eg, gehanadel@med.asu.edu.eg"
$ws.Range("G23").Value = "heba@med.asu.edu.eg, hananragab@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg"
$ws.Range("G24").Value = "NadaMohamed@med.asu.edu.eg, Amera.a.saad@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg, AbeerRagheb@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg"
$ws.Range("G25").Value = "norhan.mohamed@med.asu.edu.eg, yasmintarek@med.asu.edu.eg"
$ws.Range("G28").Value = "nourhan.osama@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, dina.adel@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, Sarah.Abdelmohsen@med.asu.edu.eg"
$ws.Range("G29").Value = "amira.m.ibrahim@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg"
$ws.Range("G31").Value = "dalia.t.abdelmegid@med.asu.edu.eg, heba.e.mohamed@med.asu.edu.eg, mona.I.hussein@med.asu.edu.eg"
$ws.Range("G32").Value = "Noran.Mahmoud@med.asu.edu.eg, menna-allah.gamil@med.asu.edu.eg"
$ws.Range("G36").Value = "enas.omran@med.asu.edu.eg, user@user.com, Walaa.h.ghanima@med.asu.edu.eg"
$ws.Range("G37").Value = "abdullah.elagrody@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, eman.samir@med.asu.edu.eg"
$ws.Range("G38").Value = "marinasorial@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, marina_atef@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg"
$ws.Range("G39").Value = "marinasorial@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, marina_atef@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg"
$ws.Range("G40").Value = "heba@med.asu.edu.eg, servinaz@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm, Amira.Sobhy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg"
$ws.Range("G41").Value = "Amira.Sobhy@med.asu.edu.eg, hananragab@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm, Mohammedeltanany@med.asu.edu.eg"
$ws.Range("G42").Value = "eman.tantawi@med.asu.edu.eg, servinaz@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm"
$ws.Range("G43").Value = "NadaMohamed@med.asu.edu.eg, Kerelos.zareef@med.asu.edu.eg, Amera.a.saad@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg, AbeerRagheb@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg, Fatmaelhady@med.asu.edu.eg"
$ws.Range("G47").Value = "nourhan.osama@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, merna.said@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg"
$ws.Range("G48").Value = "Eman.m.abosakaya@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, merna.said@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, Sarah.Abdelmohsen@med.asu.edu.eg"
$ws.Range("G50").Value = "dalia.t.abdelmegid@med.asu.edu.eg, heba.e.mohamed@med.asu.edu.eg, mona.I.hussein@med.asu.edu.eg"
$ws.Range("G51").Value = "Noran.Mahmoud@med.asu.edu.eg, menna-allah.gamil@med.asu.edu.eg"
$ws.Range("G56").Value = "abdullah.elagrody@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, eman.samir@med.asu.edu.eg"
$ws.Range("G57").Value = "marinasorial@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, marina_atef@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg"
$ws.Range("G58").Value = "marinasorial@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, marina_atef@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg"
$ws.Range("G59").Value = "heba@med.asu.edu.eg, servinaz@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg"
$ws.Range("G60").Value = "Amira.Sobhy@med.asu.edu.eg, hananragab@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm, Mohammedeltanany@med.asu.edu.eg"
$ws.Range("G63").Value = "aya.saeed@med.asu.edu.eg, aml.awwad@med.asu.edu.eg, Safa.hany@med.asu.edu.eg"
$ws.Range("G66").Value = "Marina.youhana@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, dina.adel@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg"
$ws.Range("G67").Value = "amira.m.ibrahim@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg"
$ws.Range("G71").Value = "Noran.Mahmoud@med.asu.edu.eg, Sarah.Mahdy@med.asu.edu.eg"
$ws.Range("G72").Value = "Noran.Mahmoud@med.asu.edu.eg, Sarah.Mahdy@med.asu.edu.eg"
$ws.Range("G75").Value = "abdullah.elagrody@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, eman.samir@med.asu.edu.eg"
$ws.Range("G76").Value = "marinasorial@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, eman.samir@med.asu.edu.eg"
$ws.Range("G77").Value = "marinasorial@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, marina_atef@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg"
$ws.Range("G78").Value = "heba@med.asu.edu.eg, servinaz@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg"
$ws.Range("G79").Value = "Amira.Sobhy@med.asu.edu.eg, hananragab@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm, Mohammedeltanany@med.asu.edu.eg"
$ws.Range("G81").Value = "NadaMohamed@med.asu.edu.eg, Amera.a.saad@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg, AbeerRagheb@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg"
$ws.Range("G82").Value = "norhan.mohamed@med.asu.edu.eg, yasmintarek@med.asu.edu.eg"
$ws.Range("G83").Value = "aya.saeed@med.asu.edu.eg, aml.awwad@med.asu.edu.eg, Safa.hany@med.asu.edu.eg"
$ws.Range("G85").Value = "Marina.youhana@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, dina.adel@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg"
$ws.Range("G86").Value = "Eman.m.abosakaya@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, merna.said@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, Sarah.Abdelmohsen@med.asu.edu.eg"
$ws.Range("G88").Value = "dalia.t.abdelmegid@med.asu.edu.eg, heba.e.mohamed@med.asu.edu.eg, mona.I.hussein@med.asu.edu.eg"
$ws.Range("G89").Value = "youstina.gamil@med.asu.edu.eg, Sarah.Mahdy@med.asu.edu.eg"
$ws.Range("G90").Value = "Noran.Mahmoud@med.asu.edu.eg, Sarah.Mahdy@med.asu.edu.eg"
$ws.Range("G91").Value = "Noran.Mahmoud@med.asu.edu.eg, Sarah.Mahdy@med.asu.edu.eg"
$ws.Range("G94").Value = "abdullah.elagrody@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, eman.samir@med.asu.edu.eg"
$ws.Range("G95").Value = "marinasorial@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, eman.samir@med.asu.edu.eg"
$ws.Range("G96").Value = "marinasorial@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, marina_atef@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg"
$ws.Range("G97").Value = "heba@med.asu.edu.eg, servinaz@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg"
$ws.Range("G98").Value = "Amira.Sobhy@med.asu.edu.eg, hananragab@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm, Mohammedeltanany@med.asu.edu.eg"
$ws.Range("G99").Value = "eman.tantawi@med.asu.edu.eg, servinaz@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm"
$ws.Range("G100").Value = "NadaMohamed@med.asu.edu.eg, Kerelos.zareef@med.asu.edu.eg, Amera.a.saad@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg, AbeerRagheb@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg, Fatmaelhady@med.asu.edu.eg"
$ws.Range("G101").Value = "aya.saeed@med.asu.edu.eg, aml.awwad@med.asu.edu.eg, Safa.hany@med.asu.edu.eg"
$ws.Range("G104").Value = "nourhan.osama@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, merna.said@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg"
$ws.Range("G105").Value = "amira.m.ibrahim@med.asu.edu.eg, basma.hamed@med.asu.edu.eg"
$ws.Range("G113").Value = "abdullah.elagrody@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, eman.samir@med.asu.edu.eg"
$ws.Range("G115").Value = "marinasorial@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, marina_atef@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg"
